$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value guaranteed to remain text (matches source cells, which
# are all stored as inline strings / Text type) even when the value looks
# numeric (e.g. "9", "32.0"). A leading apostrophe is Excel's standard
# "force text" entry trick and keeps the General number format.
function Set-TextValue($addr, $value) {
    $r = $ws.Range($addr)
    if ($value -match '^-?[0-9]+(\.[0-9]+)?$') {
        $r.Value = "'" + $value
    } else {
        $r.Value = $value
    }
}

# Row 6 (child 0): id 9 / Letha / Stephenie / coords / contact / (new) H6 travel time
Set-TextValue "B6" "9"
Set-TextValue "C6" "Letha  "
Set-TextValue "D6" "Stephenie  "
Set-TextValue "E6" "-9.1,7.31"
Set-TextValue "F6" "Sibyl(mother): 0567328221"
Set-TextValue "H6" "32.0"

# Row 7 (child 1)
Set-TextValue "B7" "6"
Set-TextValue "C7" "Ema  "
Set-TextValue "D7" "Ardell  "
Set-TextValue "E7" "-6.44,3.18"
Set-TextValue "F7" "Carley(grandmother): 0533587167"
Set-TextValue "G7" "7:06:00"
Set-TextValue "H7" "26.0"

# Row 8 (child 2)
Set-TextValue "B8" "8"
Set-TextValue "C8" "Marni  "
Set-TextValue "D8" "Shanika  "
Set-TextValue "E8" "-2.69,6.26"
Set-TextValue "F8" "Lady(mother): 0560804012"
Set-TextValue "G8" "7:12:00"
Set-TextValue "H8" "20.0"

# Row 9 (child 3)
Set-TextValue "B9" "3"
Set-TextValue "C9" "Alexia  "
Set-TextValue "D9" "Ramonita  "
Set-TextValue "E9" "-2.83,7.67"
Set-TextValue "F9" "Han(father): 0567537032"
Set-TextValue "G9" "7:14:00"
Set-TextValue "H9" "18.0"

# Row 10 (child 4)
Set-TextValue "B10" "19"
Set-TextValue "C10" "Jeanine  "
Set-TextValue "D10" "Janee  "
Set-TextValue "E10" "-1.93,9.03"
Set-TextValue "F10" "Teresa(mother): 0517627420"
Set-TextValue "G10" "7:16:00"
Set-TextValue "H10" "16.0"

# Row 11 (child 5)
Set-TextValue "B11" "12"
Set-TextValue "C11" "Frankie  "
Set-TextValue "D11" "Flavia  "
Set-TextValue "E11" "3.22,4.01"
Set-TextValue "F11" "Cyrus(mother): 0522363358"
Set-TextValue "G11" "7:26:00"
Set-TextValue "H11" "6.0"

# Row 12 (school) - only the arrival time changes
Set-TextValue "G12" "7:32:00"

# Row 14 (time) - total travel time
Set-TextValue "B14" "32.0"
